# Append the latest two daily data points (2025-11-27 / serial 45988 and
# 2025-11-28 / serial 45989, the latter with a zero placeholder) to every
# sheet in the workbook. Each sheet has a "date" column (A) and a
# "remn_amt" column (B); column A carries a date-time number format that
# must be carried over to the newly appended cells.

$wb = $excel.ActiveWorkbook

# Per-sheet closing balance for the new 45988 row (45989 is always 0).
$newData = @{
    "삼성바이오로직스" = 687098
    "셀트리온"         = 2236477
    "SK바이오팜"       = 380707
    "한올바이오파마"   = 155965
    "지투지바이오"     = 30908
    "대웅제약"         = 53353
}

foreach ($ws in $wb.Worksheets) {
    $name = $ws.Name
    if (-not $newData.ContainsKey($name)) {
        continue
    }

    $lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
    $dateFormat = $ws.Range("A" + $lastRow).NumberFormat

    $row1 = $lastRow + 1
    $row2 = $lastRow + 2

    $ws.Range("A" + $row1).Value = 45988
    $ws.Range("B" + $row1).Value = $newData[$name]

    $ws.Range("A" + $row2).Value = 45989
    $ws.Range("B" + $row2).Value = 0

    $ws.Range("A" + $row1).NumberFormat = $dateFormat
    $ws.Range("A" + $row2).NumberFormat = $dateFormat
}
